$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# ALC sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 54 - previously empty H:L, now populated; add new M54/N54
$ws.Range("H54").Value = 17692
$ws.Range("I54").Value = 6538
$ws.Range("J54").Value = 40000
$ws.Range("K54").Value = 6538
$ws.Range("L54").Value = 40000
$ws.Range("M54").Value = -6052
$ws.Range("N54").Value = -40972

# Row 76
$ws.Range("H76").Value = 3192.5405
$ws.Range("I76").Value = 3180.6858
$ws.Range("J76").Value = 3400
$ws.Range("K76").Value = 3180.6858
$ws.Range("L76").Value = 3400
$ws.Range("M76").Value = -2865.6858
$ws.Range("N76").Value = -4030

# Row 79
$ws.Range("H79").Value = 3192.5405
$ws.Range("I79").Value = 3180.6858
$ws.Range("J79").Value = 3400
$ws.Range("K79").Value = 3180.6858
$ws.Range("L79").Value = 3400
$ws.Range("M79").Value = -2088.6858
$ws.Range("N79").Value = -5584

# Row 132
$ws.Range("H132").Value = 131449.95
$ws.Range("I132").Value = 140517.5
$ws.Range("K132").Value = 421552.5
$ws.Range("M132").Value = -419022.5

# Row 133
$ws.Range("H133").Value = 49593.332
$ws.Range("J133").Value = 49593.332
$ws.Range("L133").Value = 49593.332
$ws.Range("N133").Value = -59713.332

# Row 134
$ws.Range("H134").Value = 59439.75
$ws.Range("J134").Value = 59439.75
$ws.Range("L134").Value = 59439.75
$ws.Range("N134").Value = -69579.75

# Row 137
$ws.Range("H137").Value = 32260024
$ws.Range("I137").Value = 1353.2916
$ws.Range("J137").Value = 142861180
$ws.Range("K137").Value = 4059.8748
$ws.Range("L137").Value = 428583540
$ws.Range("M137").Value = -1509.8748
$ws.Range("N137").Value = -428588640

# ---------------------------------------------------------------
# ARM sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 4
$ws.Range("H4").Value = 316.83334
$ws.Range("I4").Value = 300.33334
$ws.Range("J4").Value = 333.33334
$ws.Range("K4").Value = 300.33334
$ws.Range("L4").Value = 333.33334
$ws.Range("M4").Value = -184.33334
$ws.Range("N4").Value = -565.33334

# Row 32
$ws.Range("H32").Value = 4716.1113
$ws.Range("I32").Value = 4699.6787
$ws.Range("K32").Value = 4699.6787
$ws.Range("M32").Value = -4412.6787

# ---------------------------------------------------------------
# BSM sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 134
$ws.Range("H134").Value = 42293.406
$ws.Range("I134").Value = 45476.88
$ws.Range("K134").Value = 136430.64
$ws.Range("M134").Value = -133895.64

# ---------------------------------------------------------------
# CRP sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 1478.8235
$ws.Range("I31").Value = 1089.3636
$ws.Range("J31").Value = 2192.8333
$ws.Range("K31").Value = 1089.3636
$ws.Range("L31").Value = 2192.8333
$ws.Range("M31").Value = -794.3635999999999
$ws.Range("N31").Value = -2782.8333

# Row 34
$ws.Range("H34").Value = 1478.8235
$ws.Range("I34").Value = 1089.3636
$ws.Range("J34").Value = 2192.8333
$ws.Range("K34").Value = 1089.3636
$ws.Range("L34").Value = 2192.8333
$ws.Range("M34").Value = -887.3635999999999
$ws.Range("N34").Value = -2596.8333

# Row 38 - new M38 added
$ws.Range("H38").Value = 7631.143
$ws.Range("I38").Value = 1769
$ws.Range("J38").Value = 9976
$ws.Range("K38").Value = 1769
$ws.Range("L38").Value = 9976
$ws.Range("M38").Value = -1392
$ws.Range("N38").Value = -10730

# Row 46 - new M46 added
$ws.Range("H46").Value = 7631.143
$ws.Range("I46").Value = 1769
$ws.Range("J46").Value = 9976
$ws.Range("K46").Value = 1769
$ws.Range("L46").Value = 9976
$ws.Range("M46").Value = -1558
$ws.Range("N46").Value = -10398

# Row 58
$ws.Range("H58").Value = 3019.5557
$ws.Range("I58").Value = 1319
$ws.Range("K58").Value = 1319
$ws.Range("M58").Value = -1116

# Rows 129-141: clear out stale market-data columns (H:N minus G)
# Row 129 (no M129)
$ws.Range("H129:L129").ClearContents()
$ws.Range("N129").ClearContents()

# Row 130 (no M130)
$ws.Range("H130:L130").ClearContents()
$ws.Range("N130").ClearContents()

# Row 131 (no M131/N131)
$ws.Range("H131:L131").ClearContents()

# Row 132
$ws.Range("H132:L132").ClearContents()
$ws.Range("M132:N132").ClearContents()

# Row 133
$ws.Range("H133:L133").ClearContents()
$ws.Range("M133:N133").ClearContents()

# Row 134
$ws.Range("H134:L134").ClearContents()
$ws.Range("M134:N134").ClearContents()

# Row 135 (no M135)
$ws.Range("H135:L135").ClearContents()
$ws.Range("N135").ClearContents()

# Row 136 - updated, not cleared
$ws.Range("H136").Value = 3019.5557
$ws.Range("I136").Value = 1319
$ws.Range("K136").Value = 3957
$ws.Range("M136").Value = -1407

# Row 137 (no M137)
$ws.Range("H137:L137").ClearContents()
$ws.Range("N137").ClearContents()

# Row 138
$ws.Range("H138:L138").ClearContents()
$ws.Range("M138:N138").ClearContents()

# Row 139 (no M139)
$ws.Range("H139:L139").ClearContents()
$ws.Range("N139").ClearContents()

# Row 140
$ws.Range("H140:L140").ClearContents()
$ws.Range("M140:N140").ClearContents()

# Row 141 (no M141)
$ws.Range("H141:L141").ClearContents()
$ws.Range("N141").ClearContents()

# ---------------------------------------------------------------
# CUL sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 57 - N57 removed entirely
$ws.Range("H57").Value = 900
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 113
$ws.Range("H113").Value = 557.1539
$ws.Range("I113").Value = 604
$ws.Range("J113").Value = 527.875
$ws.Range("K113").Value = 1812
$ws.Range("L113").Value = 1583.625
$ws.Range("M113").Value = 358
$ws.Range("N113").Value = -5923.625

# ---------------------------------------------------------------
# GSM sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 132
$ws.Range("H132").Value = 1988.3243
$ws.Range("I132").Value = 1902.5
$ws.Range("J132").Value = 2146.7693
$ws.Range("K132").Value = 5707.5
$ws.Range("L132").Value = 6440.3079
$ws.Range("M132").Value = -3177.5
$ws.Range("N132").Value = -11500.3079

# ---------------------------------------------------------------
# LTW sheet
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 132
$ws.Range("H132").Value = 4089.4187
$ws.Range("I132").Value = 5509.636
$ws.Range("J132").Value = 2601.5715
$ws.Range("K132").Value = 16528.908
$ws.Range("L132").Value = 7804.7145
$ws.Range("M132").Value = -13998.908
$ws.Range("N132").Value = -12864.7145
